# Example Dictionary.xlsx - add "ID" variable rows to both the ETC and RED
# survey blocks.
#
# 1. Insert a new row after the existing ETC block (row 11) that documents
#    an "ID" variable, pushing the RED block down by one row.
# 2. Append a matching "ID" row at the end of the RED block (new row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new ETC "ID" row right after the ETC/UNLINKED_WGHT_FCTR row,
# before the RED block begins.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "ETC"
$ws.Range("B11").Value = "ID"
$ws.Range("C11").Value = "NONCATEGORICAL"
$ws.Range("D11").Value = "ID"
$ws.Range("E11").Value = "NONCATEGORICAL"

# Append the new RED "ID" row at the end of the table.
$ws.Range("A20").Value = "RED"
$ws.Range("B20").Value = "ID"
$ws.Range("C20").Value = "NONCATEGORICAL"
$ws.Range("D20").Value = "ID"
$ws.Range("E20").Value = "NONCATEGORICAL"

# Leave the selection where the author's saved workbook shows it.
$ws.Range("A21").Select()
